$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the four evaluation comments whose wording was tweaked ---
# G2 (Task 1.1 reason for "Ababei Vasile")
$ws.Range("G2").Value = ' The function respects the requirement. It correctly implements a BFS algorithm to find the shortest path between the start and end positions in the matrix, considering the given positions as obstacles. The function returns the minimum path length, which is 4 in this case, as expected.'

# I2 (Task 1.2 reason for "Ababei Vasile")
$ws.Range("I2").Value = ' The function respects the requirement because it performs a BFS traversal of the matrix, which has a time complexity of O(n * m) in the worst case, and it also iterates over the given positions, which has a time complexity of O(k). Therefore, the overall time complexity is O(n * m + k), which is at most O(n * m * k) since k is less than or equal to n * m.'

# K2 (Task 1.3 reason for "Ababei Vasile")
$ws.Range("K2").Value = ' The function respects the requirement because it uses a visited array of size n*m and a queue of maximum size n*m, resulting in a space complexity of O(n*m).'

# G3 (Task 1.1 reason for "Alexe Robert George")
$ws.Range("G3").Value = ' The function respects the requirement. It correctly implements a BFS algorithm to find the shortest path from the starting position to the ending position while avoiding the given positions. The function returns the minimum path length, which is 4 in this case, as expected.'

# --- Columns G and I narrow slightly to fit the new (shorter) longest text in each column ---
$ws.Columns.Item(7).ColumnWidth = 248.66666666666666
$ws.Columns.Item(9).ColumnWidth = 306.8333333333333
